$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.572.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.45"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +13.96%  "
$ws.Range("E9").Value = "  +8.21%  "
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.111.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.43"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.37%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.675"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.24%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.831.96"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.567.90"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.58%  "
$ws.Range("E22").Value = "  +14.00%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +31.87%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.352.50"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +37.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0566"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +9.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.96"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +16.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.697"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.14"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.345.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("E40").Value = "  +6.02%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.55%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.43"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("E43").Value = "  +7.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.27"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.27"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0514"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.018.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.69%  "
